$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$E32 = '{"0": [{"frameNodeType": "MotionFrameNode_Animation", "frameGuid": "6bdb", "guid": "29769", "duration": 12, "animSlot": -1}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "e1ca", "guid": "145503", "offsetPos": {"x": 0, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": -12, "slotIndex": 17, "colorHex": "", "colorHex1": ""}], "1": [], "2": [{"frameNodeType": "MotionFrameNode_Sound", "frameGuid": "be8a", "guid": "137576", "sound_volume": 1, "sound_innerRadius": 100, "sound_maxDistance": 3000}], "3": [{"frameNodeType": "MotionFrameNode_impulse", "frameGuid": "30b5", "offsetPos": {"x": -1000, "y": 0, "z": 0}, "forceNum": 3000, "groundFriction": 1, "gravityScale": 1, "gravityTime": 0}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "ac13", "guid": "151743", "offsetPos": {"x": 0, "y": 0, "z": -70}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 4, "y": 4, "z": 4}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}], "4": [{"frameNodeType": "MotionFrameNode_Shake", "frameGuid": "23e8", "rollAmplitude": 0, "rollFrequency": 0, "rollWaveform": 1, "pitchAmplitude": 0, "pitchFrequency": 0, "pitchWaveform": 1, "yawAmplitude": 0, "yawFrequency": 0, "yawWaveform": 1, "xAmplitude": 6, "xFrequency": 10, "xWaveform": 1, "yAmplitude": 0, "yFrequency": 0, "yWaveform": 1, "zAmplitude": 6, "zFrequency": 4, "zWaveform": 1, "fovAmplitude": 0, "fovFrequency": 0, "fovWaveform": 1, "keepTime": 0.4}, {"frameNodeType": "MotionFrameNode_FlyEntity", "frameGuid": "9c5a", "delayTime": 0, "bulletId": 19, "startLoc": {"x": 100, "y": 0, "z": 30}, "endType": 0, "endDis": 2000}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "fffb", "guid": "291299", "offsetPos": {"x": 40, "y": 0, "z": 50}, "offsetRotation": {"x": 0, "y": 0, "z": 90}, "offsetScale": {"x": 1.5, "y": 1.5, "z": 1.5}, "count": -9, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#DDD300FF"}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "f3ae", "type": 0, "offsetLoc": {"x": 30, "y": 0, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 5, "y": 10, "z": 4}, "range": 1, "checkCount": 3, "checkInterval": 0.1, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 26, "dilationRate": 0, "dilationFrame": 0}], "5": [{"frameNodeType": "MotionFrameNode_AnimPause", "frameGuid": "476f", "guid": "29769", "pauseCount": 2}], "6": [], "7": [{"frameNodeType": "MotionFrameNode_BreakPoint", "frameGuid": "568f", "breakType": 2}], "8": []}'
$E33 = '{"0": [{"frameNodeType": "MotionFrameNode_Animation", "frameGuid": "b279", "guid": "29747", "duration": 10, "animSlot": 0}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "cd37", "guid": "27694", "offsetPos": {"x": 0, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 3, "y": 3, "z": 3}, "count": -7, "slotIndex": 12, "colorHex": "", "colorHex1": "Color|#B200FFFF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "fcbe", "guid": "27694", "offsetPos": {"x": 0, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 3, "y": 3, "z": 3}, "count": -7, "slotIndex": 12, "colorHex": "", "colorHex1": "Color|#B200FFFF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "c629", "guid": "163346", "offsetPos": {"x": 0, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": -2, "slotIndex": 16, "colorHex": "", "colorHex1": ""}], "1": [{"frameNodeType": "MotionFrameNode_AnimPause", "frameGuid": "62d7", "guid": "29747", "pauseCount": 1}], "2": [], "3": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "a19b", "guid": "295654", "offsetPos": {"x": 500, "y": 0, "z": -50}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "33dd", "type": 1, "offsetLoc": {"x": 500, "y": 0, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 10, "checkCount": 1, "checkInterval": 0, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 27, "dilationRate": 0, "dilationFrame": 0}, {"frameNodeType": "MotionFrameNode_3DSound", "frameGuid": "dbb8", "guid": "137576", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000, "sound_count": 0, "bind": 1}, {"frameNodeType": "MotionFrameNode_Shake", "frameGuid": "7ad3", "rollAmplitude": 0, "rollFrequency": 0, "rollWaveform": 1, "pitchAmplitude": 0, "pitchFrequency": 0, "pitchWaveform": 1, "yawAmplitude": 0, "yawFrequency": 0, "yawWaveform": 1, "xAmplitude": 20, "xFrequency": 5, "xWaveform": 0, "yAmplitude": 10, "yFrequency": 2, "yWaveform": 1, "zAmplitude": 14, "zFrequency": 7, "zWaveform": 1, "fovAmplitude": 0, "fovFrequency": 0, "fovWaveform": 1, "keepTime": 0.2}], "4": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "fa6f", "guid": "295654", "offsetPos": {"x": 900, "y": 0, "z": -50}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "d7d6", "type": 1, "offsetLoc": {"x": 900, "y": 0, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 10, "checkCount": 1, "checkInterval": 0, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 27, "dilationRate": 0, "dilationFrame": 0}, {"frameNodeType": "MotionFrameNode_3DSound", "frameGuid": "716f", "guid": "137576", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000, "sound_count": 0, "bind": 1}], "5": [{"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "735b", "type": 1, "offsetLoc": {"x": 1300, "y": 0, "z": -80}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 10, "checkCount": 1, "checkInterval": 0, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 27, "dilationRate": 0, "dilationFrame": 0}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "9b01", "guid": "295654", "offsetPos": {"x": 1300, "y": 0, "z": -80}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_3DSound", "frameGuid": "5a34", "guid": "137576", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000, "sound_count": 0, "bind": 1}], "6": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "e6a2", "guid": "295654", "offsetPos": {"x": 1700, "y": 0, "z": -80}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_BreakPoint", "frameGuid": "44f6", "breakType": 2}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "b99a", "type": 1, "offsetLoc": {"x": 1700, "y": 0, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 10, "checkCount": 1, "checkInterval": 0, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 27, "dilationRate": 0, "dilationFrame": 0}, {"frameNodeType": "MotionFrameNode_3DSound", "frameGuid": "3d14", "guid": "137576", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000, "sound_count": 0, "bind": 1}], "7": [], "8": [], "9": [], "11": [], "12": []}'
$E35 = '{"0": [{"frameNodeType": "MotionFrameNode_Animation", "frameGuid": "90fd", "guid": "122510", "duration": 30, "animSlot": 0}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "42ec", "guid": "27694", "offsetPos": {"x": 0, "y": 0, "z": 100}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 6, "y": 6, "z": 6}, "count": -10, "slotIndex": 23, "colorHex": "", "colorHex1": "Color|#0089FFFF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "88df", "guid": "27694", "offsetPos": {"x": 0, "y": 0, "z": 100}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 6, "y": 6, "z": 6}, "count": -10, "slotIndex": 23, "colorHex": "", "colorHex1": "Color|#0089FFFF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "0c41", "guid": "181010", "offsetPos": {"x": 0, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1.2, "y": 1.2, "z": 1.2}, "count": -10, "slotIndex": 16, "colorHex": "", "colorHex1": "Color|#0053FFFF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "2dfd", "guid": "168954", "offsetPos": {"x": 0, "y": 0, "z": -70}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 2, "y": 2, "z": 4}, "count": -11, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#6F92B8FF"}, {"frameNodeType": "MotionFrameNode_Shake", "frameGuid": "7f11", "rollAmplitude": 0, "rollFrequency": 0, "rollWaveform": 1, "pitchAmplitude": 0, "pitchFrequency": 0, "pitchWaveform": 1, "yawAmplitude": 0, "yawFrequency": 0, "yawWaveform": 1, "xAmplitude": 4, "xFrequency": 4, "xWaveform": 1, "yAmplitude": 0, "yFrequency": 0, "yWaveform": 1, "zAmplitude": 4, "zFrequency": 4, "zWaveform": 1, "fovAmplitude": 0, "fovFrequency": 0, "fovWaveform": 1, "keepTime": 1}], "5": [{"frameNodeType": "MotionFrameNode_Sound", "frameGuid": "c4f4", "guid": "265134", "sound_volume": 2, "sound_innerRadius": 20000, "sound_maxDistance": 3000}], "7": [], "11": [{"frameNodeType": "MotionFrameNode_Shake", "frameGuid": "d11a", "rollAmplitude": 0, "rollFrequency": 0, "rollWaveform": 1, "pitchAmplitude": 0, "pitchFrequency": 0, "pitchWaveform": 1, "yawAmplitude": 0, "yawFrequency": 0, "yawWaveform": 1, "xAmplitude": 10, "xFrequency": 10, "xWaveform": 1, "yAmplitude": 0, "yFrequency": 0, "yWaveform": 1, "zAmplitude": 10, "zFrequency": 10, "zWaveform": 1, "fovAmplitude": 50, "fovFrequency": 0, "fovWaveform": 1, "keepTime": 0.5}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "20d0", "guid": "89107", "offsetPos": {"x": 150, "y": 0, "z": -70}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 2, "y": 2, "z": 2}, "count": -6, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#009BFFFF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "bcac", "guid": "89107", "offsetPos": {"x": 150, "y": 0, "z": -70}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": -6, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#009BFFFF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "d85a", "guid": "168949", "offsetPos": {"x": 150, "y": 0, "z": -70}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 3, "y": 3, "z": 3}, "count": -5, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#066EFFFF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "19eb", "guid": "135892", "offsetPos": {"x": 150, "y": 0, "z": -70}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 3, "y": 3, "z": 3}, "count": -5, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#066EFFFF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "02e8", "guid": "130642", "offsetPos": {"x": 150, "y": 0, "z": -70}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 7, "y": 7, "z": 7}, "count": -8, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#066EFFFF"}], "12": [{"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "c5f0", "type": 1, "offsetLoc": {"x": 150, "y": 0, "z": -80}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 11, "checkCount": 5, "checkInterval": 0.1, "keepFrameCount": 5, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 31, "dilationRate": 0, "dilationFrame": 0}], "14": [{"frameNodeType": "MotionFrameNode_BreakPoint", "frameGuid": "7e7a", "breakType": 1}], "15": [], "16": [], "18": []}'
$E36 = '{"0": [{"frameNodeType": "MotionFrameNode_Animation", "frameGuid": "d536", "guid": "279656", "duration": 16, "animSlot": 0}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "8a78", "type": 0, "offsetLoc": {"x": 0, "y": 0, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 1, "checkCount": 1, "checkInterval": 0, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 46, "dilationRate": 0, "dilationFrame": 0}], "1": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "d90c", "guid": "158173", "offsetPos": {"x": 0, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 0.5, "y": 0.5, "z": 0.5}, "count": -3, "slotIndex": -1, "colorHex": "", "colorHex1": ""}], "2": [{"frameNodeType": "MotionFrameNode_AnimPause", "frameGuid": "b2a2", "guid": "279656", "pauseCount": 1}], "3": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "236a", "guid": "295658", "offsetPos": {"x": 0, "y": 0, "z": -20}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_Shake", "frameGuid": "9f35", "rollAmplitude": 0, "rollFrequency": 0, "rollWaveform": 1, "pitchAmplitude": 0, "pitchFrequency": 0, "pitchWaveform": 1, "yawAmplitude": 0, "yawFrequency": 0, "yawWaveform": 1, "xAmplitude": 20, "xFrequency": 0, "xWaveform": 1, "yAmplitude": 0, "yFrequency": 0, "yWaveform": 1, "zAmplitude": 16, "zFrequency": 3, "zWaveform": 1, "fovAmplitude": 100, "fovFrequency": 10, "fovWaveform": 1, "keepTime": 0.5}], "4": [{"frameNodeType": "MotionFrameNode_Move", "frameGuid": "f36b", "move_offsetPos": {"x": 1000, "y": 0, "z": 0}, "move_isFlash": 0, "move_count": 0}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "c2cb", "guid": "200145", "offsetPos": {"x": -1400, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 0.5, "y": 0.7, "z": 0.7}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_Sound", "frameGuid": "9255", "guid": "137559", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "de14", "type": 0, "offsetLoc": {"x": -400, "y": 0, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 15, "y": 4, "z": 4}, "range": 1, "checkCount": 1, "checkInterval": 0.1, "keepFrameCount": 2, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 32, "dilationRate": 0, "dilationFrame": 0}], "6": [{"frameNodeType": "MotionFrameNode_BreakPoint", "frameGuid": "15b4", "breakType": 1}], "7": [], "8": [], "9": [], "10": [], "12": [], "13": [], "14": [], "15": [], "16": [], "17": []}'
$E37 = '{"0": [{"frameNodeType": "MotionFrameNode_Animation", "frameGuid": "ad79", "guid": "285460", "duration": 16, "animSlot": 0}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "aa84", "guid": "163346", "offsetPos": {"x": 0, "y": 0, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": -6, "slotIndex": 16, "colorHex": "", "colorHex1": ""}], "5": [{"frameNodeType": "MotionFrameNode_AnimPause", "frameGuid": "2b26", "guid": "285460", "pauseCount": 1}], "6": [], "7": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "de5d", "guid": "295655", "offsetPos": {"x": 240, "y": 40, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": "Color#|AE4448FF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "d848", "guid": "295655", "offsetPos": {"x": -480, "y": 70, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": "Color#|AE4448FF"}, {"frameNodeType": "MotionFrameNode_3DSound", "frameGuid": "70b8", "guid": "137576", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000, "sound_count": 0, "bind": 1}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "4e3c", "type": 1, "offsetLoc": {"x": 240, "y": 40, "z": 60}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 6, "checkCount": 1, "checkInterval": 0.2, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 32, "dilationRate": 0, "dilationFrame": 0}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "35ea", "type": 1, "offsetLoc": {"x": -480, "y": 70, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 6, "checkCount": 1, "checkInterval": 0.2, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 32, "dilationRate": 0, "dilationFrame": 0}], "8": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "7afa", "guid": "295655", "offsetPos": {"x": 60, "y": 320, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#AE4448FF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "dfdd", "guid": "295655", "offsetPos": {"x": 120, "y": -480, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#AE4448FF"}, {"frameNodeType": "MotionFrameNode_3DSound", "frameGuid": "85b9", "guid": "137576", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000, "sound_count": 0, "bind": 1}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "9f79", "type": 1, "offsetLoc": {"x": 60, "y": 320, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 6, "checkCount": 1, "checkInterval": 0.2, "keepFrameCount": 1, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 32, "dilationRate": 0, "dilationFrame": 0}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "af0a", "type": 1, "offsetLoc": {"x": 120, "y": -480, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 6, "checkCount": 1, "checkInterval": 0.2, "keepFrameCount": 1, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 32, "dilationRate": 0, "dilationFrame": 0}], "9": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "f2b7", "guid": "295655", "offsetPos": {"x": 590, "y": 430, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#AE4448FF"}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "2223", "guid": "295655", "offsetPos": {"x": -670, "y": -410, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": "Color|#AE4448FF"}, {"frameNodeType": "MotionFrameNode_3DSound", "frameGuid": "5988", "guid": "137576", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000, "sound_count": 0, "bind": 1}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "c583", "type": 1, "offsetLoc": {"x": 590, "y": 430, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 6, "checkCount": 1, "checkInterval": 0.2, "keepFrameCount": 1, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 32, "dilationRate": 0, "dilationFrame": 0}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "4e47", "type": 1, "offsetLoc": {"x": -670, "y": -410, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 6, "checkCount": 1, "checkInterval": 0.2, "keepFrameCount": 1, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 32, "dilationRate": 0, "dilationFrame": 0}], "10": [{"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "ecdd", "guid": "295655", "offsetPos": {"x": -420, "y": 370, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_Effect", "frameGuid": "1476", "guid": "295655", "offsetPos": {"x": 590, "y": 130, "z": 0}, "offsetRotation": {"x": 0, "y": 0, "z": 0}, "offsetScale": {"x": 1, "y": 1, "z": 1}, "count": 1, "slotIndex": -1, "colorHex": "", "colorHex1": ""}, {"frameNodeType": "MotionFrameNode_3DSound", "frameGuid": "d1eb", "guid": "137576", "sound_volume": 10, "sound_innerRadius": 100, "sound_maxDistance": 3000, "sound_count": 0, "bind": 1}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "2ea1", "type": 1, "offsetLoc": {"x": 590, "y": 370, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 6, "checkCount": 1, "checkInterval": 0.2, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 32, "dilationRate": 0, "dilationFrame": 0}, {"frameNodeType": "MotionFrameNode_SkillRect", "frameGuid": "8fe6", "type": 1, "offsetLoc": {"x": -420, "y": 370, "z": 0}, "offsetRot": {"x": 0, "y": 0, "z": 0}, "LWH": {"x": 1, "y": 1, "z": 1}, "range": 6, "checkCount": 1, "checkInterval": 0.2, "keepFrameCount": 0, "rectSocket": 0, "isShake": 0, "cameraShakeId": 0, "effectid": 32, "dilationRate": 0, "dilationFrame": 0}], "15": [], "16": [{"frameNodeType": "MotionFrameNode_BreakPoint", "frameGuid": "8418", "breakType": 2}], "18": [], "20": [], "21": [], "22": [], "23": [], "25": [], "26": [], "27": [], "28": []}'

$ws.Range("E32").Value = $E32
$ws.Range("E33").Value = $E33
$ws.Range("D37").Value = 29
$ws.Range("E35").Value = $E35
$ws.Range("E36").Value = $E36
$ws.Range("E37").Value = $E37

$wb.Save()
